$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.369.28"
$ws.Range("E2").Value = "  +8.86%  "
$ws.Range("D3").Value = "'1.607.92"
$ws.Range("E3").Value = "  +9.29%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'0.9941"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "'290.01"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("D7").Value = "'0.3720"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.3383"
$ws.Range("E8").Value = "  +10.69%  "
$ws.Range("D9").Value = "'42.75"
$ws.Range("E9").Value = "  +7.38%  "
$ws.Range("D10").Value = "'1.150"
$ws.Range("E10").Value = "  +8.63%  "
$ws.Range("D11").Value = "'0.07097"
$ws.Range("E11").Value = "  +6.94%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'19.96"
$ws.Range("E13").Value = "  +10.69%  "
$ws.Range("D14").Value = "'5.952"
$ws.Range("E14").Value = "  +8.81%  "
$ws.Range("D15").Value = "'6.685"
$ws.Range("E15").Value = "  +8.06%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001088"
$ws.Range("E16").Value = "  +5.80%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'1.605.66"
$ws.Range("E17").Value = "  +9.18%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9937"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "'0.06648"
$ws.Range("E19").Value = "  +12.12%  "
$ws.Range("D20").Value = "'78.61"
$ws.Range("E20").Value = "  +13.38%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'16.25"
$ws.Range("E21").Value = "  +12.34%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'6.052"
$ws.Range("E22").Value = "  +11.30%  "
$ws.Range("D23").Value = "'11.83"
$ws.Range("E23").Value = "  +7.50%  "
$ws.Range("D24").Value = "'22.400.38"
$ws.Range("E24").Value = "  +8.98%  "
$ws.Range("D25").Value = "'2.401"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("D26").Value = "'2.511"
$ws.Range("E26").Value = "  +18.60%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'150.88"
$ws.Range("E27").Value = "  +7.30%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.56"
$ws.Range("E28").Value = "  +13.84%  "
$ws.Range("D29").Value = "'1.784.55"
$ws.Range("E29").Value = "  +9.59%  "
$ws.Range("D30").Value = "'121.25"
$ws.Range("E30").Value = "  +6.41%  "
$ws.Range("D31").Value = "'4.227"
$ws.Range("E31").Value = "  +7.18%  "
$ws.Range("D32").Value = "'6.013"
$ws.Range("E32").Value = "  +21.59%  "
$ws.Range("D33").Value = "'0.9570"
$ws.Range("E33").Value = "  +18.39%  "
$ws.Range("D34").Value = "'0.08270"
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("D35").Value = "'1.633"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("D36").Value = "'5.348"
$ws.Range("E36").Value = "  +14.10%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'8.708"
$ws.Range("E37").Value = "  +14.02%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'11.85"
$ws.Range("E38").Value = "  +14.39%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06202"
$ws.Range("E39").Value = "  +6.99%  "
$ws.Range("D40").Value = "'1.235"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "'0.02211"
$ws.Range("E41").Value = "  +8.46%  "
$ws.Range("D42").Value = "'0.2031"
$ws.Range("E42").Value = "  +8.31%  "
$ws.Range("D43").Value = "'0.6008"
$ws.Range("E43").Value = "  +13.91%  "
$ws.Range("D44").Value = "'0.9929"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").Value = "'13.25"
$ws.Range("E45").Value = "  +9.55%  "
$ws.Range("D46").Value = "'3.666"
$ws.Range("E46").Value = "  +4.37%  "
$ws.Range("D47").Value = "'0.5760"
$ws.Range("E47").Value = "  +11.31%  "
$ws.Range("D48").Value = "'125.60"
$ws.Range("E48").Value = "  +5.34%  "
$ws.Range("D49").Value = "'1.982"
$ws.Range("E49").Value = "  +11.04%  "
$ws.Range("D50").Value = "'0.06870"
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("D51").Value = "'73.81"
$ws.Range("E51").Value = "  +9.94%  "
